$wb = $excel.ActiveWorkbook

# --- Logs sheet: append row 3 ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A3").Value = "Opvolging bestelling"
$logs.Range("B3").Value = "mailmind.test@zohomail.eu"
$logs.Range("D3").Value = "Inkoop / Bestellingen"
$logs.Range("F3").Value = "2025-08-30 18:34:31"
$logs.Range("G3").Value = "Nee"
$logs.Range("H3").Value = "Ja"
$logs.Range("I3").Value = "Nee"
$logs.Range("J3").Value = "Nee"

# Extend the conditional formatting ranges so they cover the new row too
$logs.Range("D2").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D3"))
$logs.Range("G2").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G3"))
$logs.Range("H2").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H3"))
$logs.Range("I2").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I3"))
$logs.Range("J2").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J3"))

# --- Dashboard sheet: update the aggregate count ---
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 2
